$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content in column H (order matters for sharedStrings indexing) ---
$ws.Range("H6").Value = "ในส่วนนี้ตัวเลขที่ทำกากรอกจำนวนไปค้างเป็ฯราคาเดิม"
$ws.Range("H4").Value = "ระบบไม่ได้ทำการดักalertในส่วนของการกรอกจำนวนไม่ถูกต้อง"
$ws.Range("H3").Value = "ในส่วนนี้ถ้ากรอกตัวเลขต้องเป็ฯจำนวนเป็ฯบวกเท่านั้นจะต้องไม่สารมารถกรอกตัวเลขติดลบไม่ได้"

# --- Merge H4:H5 and apply border + alignment formatting ---
$ws.Range("H4:H5").Merge()
$ws.Range("H4:H5").HorizontalAlignment = -4131
$ws.Range("H4:H5").VerticalAlignment = -4108
$ws.Range("H4:H5").Borders.Item(7).LineStyle = 1

# --- Column H width ---
$ws.Columns("H").ColumnWidth = 63.8658854166667

# --- Update selection to match the new active cell ---
$ws.Range("H13").Select()
